# Rebuild the tail of the Turbidity data table (rows 198-272 on Sheet1):
#  - Rows 198-214 keep their original Date/Turbidity pairs, but the Date
#    column switches from text to a real number (t="inlineStr" -> numeric).
#  - Rows 215-259 are brand new numeric-date rows appended after them.
#  - Rows 260-272 are a new block of rows whose Date column is still text
#    (quote-prefixed so it round-trips as a string, like the original
#    inlineStr cells) paired with updated Turbidity numbers.
# The sheet's used range / <dimension> grows from A1:B214 to A1:B272
# automatically once these cells are populated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 198-259, column A numeric dates, column B numeric values ---
$numData = New-Object "object[,]" 62,2
$numData[0,0] = 20230203
$numData[0,1] = 1.271000057458878
$numData[1,0] = 20230208
$numData[1,1] = 1.316999942064285
$numData[2,0] = 20230213
$numData[2,1] = 1.324999928474426
$numData[3,0] = 20230218
$numData[3,1] = 1.309999972581863
$numData[4,0] = 20230223
$numData[4,1] = 1.313000023365021
$numData[5,0] = 20230228
$numData[5,1] = 1.574999988079071
$numData[6,0] = 20230305
$numData[6,1] = 1.771000027656555
$numData[7,0] = 20230310
$numData[7,1] = 1.615999937057495
$numData[8,0] = 20230315
$numData[8,1] = 3.242000043392181
$numData[9,0] = 20230320
$numData[9,1] = 1.307000070810318
$numData[10,0] = 20230325
$numData[10,1] = 1.301999986171722
$numData[11,0] = 20230330
$numData[11,1] = 1.342000067234039
$numData[12,0] = 20230404
$numData[12,1] = 1.456999927759171
$numData[13,0] = 20230409
$numData[13,1] = 1.335999965667725
$numData[14,0] = 20230414
$numData[14,1] = 1.270000040531158
$numData[15,0] = 20230419
$numData[15,1] = 1.488000005483627
$numData[16,0] = 20230424
$numData[16,1] = 1.666000038385391
$numData[17,0] = 20230305
$numData[17,1] = 1.780000030994415
$numData[18,0] = 20230310
$numData[18,1] = 1.573999971151352
$numData[19,0] = 20230315
$numData[19,1] = 1.544999927282333
$numData[20,0] = 20230320
$numData[20,1] = 1.334999948740005
$numData[21,0] = 20230325
$numData[21,1] = 1.272000074386597
$numData[22,0] = 20230330
$numData[22,1] = 1.389999985694885
$numData[23,0] = 20230404
$numData[23,1] = 1.402000039815903
$numData[24,0] = 20230409
$numData[24,1] = 1.372999995946884
$numData[25,0] = 20230414
$numData[25,1] = 1.272999942302704
$numData[26,0] = 20230419
$numData[26,1] = 1.500999927520752
$numData[27,0] = 20230424
$numData[27,1] = 1.405999958515167
$numData[28,0] = 20230429
$numData[28,1] = 1.928000003099442
$numData[29,0] = 20230504
$numData[29,1] = 1.22299998998642
$numData[30,0] = 20230305
$numData[30,1] = 1.790000051259995
$numData[31,0] = 20230310
$numData[31,1] = 1.60300001502037
$numData[32,0] = 20230315
$numData[32,1] = 1.198000013828278
$numData[33,0] = 20230320
$numData[33,1] = 1.319999992847443
$numData[34,0] = 20230325
$numData[34,1] = 1.269000023603439
$numData[35,0] = 20230330
$numData[35,1] = 1.331000030040741
$numData[36,0] = 20230404
$numData[36,1] = 1.397999972105026
$numData[37,0] = 20230409
$numData[37,1] = 1.442999988794327
$numData[38,0] = 20230414
$numData[38,1] = 1.308999955654144
$numData[39,0] = 20230419
$numData[39,1] = 1.474000066518784
$numData[40,0] = 20230424
$numData[40,1] = 1.430000066757202
$numData[41,0] = 20230429
$numData[41,1] = 1.668000072240829
$numData[42,0] = 20230504
$numData[42,1] = 1.234999969601631
$numData[43,0] = 20230305
$numData[43,1] = 1.816000044345856
$numData[44,0] = 20230310
$numData[44,1] = 1.553999930620193
$numData[45,0] = 20230315
$numData[45,1] = 1.407999992370605
$numData[46,0] = 20230320
$numData[46,1] = 1.337999999523163
$numData[47,0] = 20230325
$numData[47,1] = 1.180000007152557
$numData[48,0] = 20230330
$numData[48,1] = 1.221000030636787
$numData[49,0] = 20230404
$numData[49,1] = 1.292999982833862
$numData[50,0] = 20230409
$numData[50,1] = 1.381999999284744
$numData[51,0] = 20230414
$numData[51,1] = 1.289999932050705
$numData[52,0] = 20230419
$numData[52,1] = 1.381999999284744
$numData[53,0] = 20230424
$numData[53,1] = 1.29600003361702
$numData[54,0] = 20230429
$numData[54,1] = 6.155999898910522
$numData[55,0] = 20230504
$numData[55,1] = 1.142999976873398
$numData[56,0] = 20220608
$numData[56,1] = 8.059999942779541
$numData[57,0] = 20220613
$numData[57,1] = 8.672000169754028
$numData[58,0] = 20220618
$numData[58,1] = 1.951999962329865
$numData[59,0] = 20220703
$numData[59,1] = 8.880000114440918
$numData[60,0] = 20220901
$numData[60,1] = 1.578000038862228
$numData[61,0] = 20220906
$numData[61,1] = 1.159000024199486
$ws.Range("A198:B259").Value = $numData

# --- Block 2: rows 260-272, column A text (quote-prefixed) dates, column B numeric values ---
$txtData = New-Object "object[,]" 13,2
$txtData[0,0] = "'20230305"
$txtData[0,1] = 1.781000047922134
$txtData[1,0] = "'20230310"
$txtData[1,1] = 1.596000045537949
$txtData[2,0] = "'20230315"
$txtData[2,1] = 1.757999956607819
$txtData[3,0] = "'20230320"
$txtData[3,1] = 1.31400004029274
$txtData[4,0] = "'20230325"
$txtData[4,1] = 1.274999976158142
$txtData[5,0] = "'20230330"
$txtData[5,1] = 1.425999999046326
$txtData[6,0] = "'20230404"
$txtData[6,1] = 1.421000063419342
$txtData[7,0] = "'20230409"
$txtData[7,1] = 1.333999931812286
$txtData[8,0] = "'20230414"
$txtData[8,1] = 1.27700001001358
$txtData[9,0] = "'20230419"
$txtData[9,1] = 1.480000019073486
$txtData[10,0] = "'20230424"
$txtData[10,1] = 1.440999954938889
$txtData[11,0] = "'20230429"
$txtData[11,1] = 1.867000013589859
$txtData[12,0] = "'20230504"
$txtData[12,1] = 1.24600000679493
$ws.Range("A260:B272").Value = $txtData

"done"